# Book Package: total Price Function
#
# Applies the changes described by the commit diff:
#  - admin_login_tbl: duplicate existing row 2 data into new row 5
#  - category_management_tbl: append a member/booking-style record (row 4)
#    plus several more category rows (rows 5-12)
#  - member_signup_tbl: overwrite row 2 with booking-style values and append
#    several duplicate/new member rows (rows 3, 5, 6, 7, 8)
#  - new sheet booking_contact_details_tbl with header + two data rows
#  - workbook defined names updated / added accordingly

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing text storage for anything
# that Excel's automatic type inference would otherwise turn into a real
# number or date (every cell in this workbook is a plain text/inline
# string, e.g. pincodes, phone numbers and dates are all stored as text).
function Set-TextValue($sheet, $ref, $value) {
    $cell = $sheet.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# ---------------------------------------------------------------------------
# 1. admin_login_tbl -- add row 5 (duplicate of row 2)
# ---------------------------------------------------------------------------
$adminSheet = $wb.Worksheets.Item("admin_login_tbl")
Set-TextValue $adminSheet "A5" "hetvee_sakaria"
Set-TextValue $adminSheet "B5" "Hetu@2001"
Set-TextValue $adminSheet "C5" "Hetvee Sakaria"

# ---------------------------------------------------------------------------
# 2. category_management_tbl -- add rows 4-12
# ---------------------------------------------------------------------------
$categorySheet = $wb.Worksheets.Item("category_management_tbl")

Set-TextValue $categorySheet "A4" "Ritu Agarwal"
Set-TextValue $categorySheet "B4" "2008-06-10"
Set-TextValue $categorySheet "C4" "54879672315"
Set-TextValue $categorySheet "D4" "rituagarwal@gmail.com"
Set-TextValue $categorySheet "E4" "Gujarat"
Set-TextValue $categorySheet "F4" "Gandhinagar"
Set-TextValue $categorySheet "G4" "345003"
Set-TextValue $categorySheet "H4" "Pratik Colony, Near Yash Complex _x000d_`nSG Highway, Jaipur, Rajasthan"
Set-TextValue $categorySheet "I4" "2"
Set-TextValue $categorySheet "J4" "123456"
Set-TextValue $categorySheet "K4" "pending"

Set-TextValue $categorySheet "B5" "National Packages"
Set-TextValue $categorySheet "B6" "Domestic Packages"
Set-TextValue $categorySheet "B7" "International Packages"
Set-TextValue $categorySheet "B8" "National Packages"
Set-TextValue $categorySheet "B9" "Domestic Packages"
Set-TextValue $categorySheet "B10" "International Packages"
Set-TextValue $categorySheet "B11" "National Packages"
Set-TextValue $categorySheet "B12" "Domestic Packages"

# ---------------------------------------------------------------------------
# 3. member_signup_tbl -- overwrite row 2, add rows 3, 5, 6, 7, 8
# ---------------------------------------------------------------------------
$memberSheet = $wb.Worksheets.Item("member_signup_tbl")

Set-TextValue $memberSheet "B2" "joepatrik@gmail.com"
Set-TextValue $memberSheet "C2" "35000"
Set-TextValue $memberSheet "D2" "2"
Set-TextValue $memberSheet "E2" "4"
Set-TextValue $memberSheet "F2" "No thank u"
# H2 is untouched by the edit, but re-assert it explicitly so its literal
# "_x000d_" + linebreak text content is preserved verbatim.
Set-TextValue $memberSheet "H2" "Pratik Colony, Near Yash Complex _x000d_`nSG Highway, Jaipur, Rajasthan"

function Set-JoePatrikRow($sheet, $row) {
    Set-TextValue $sheet "A$row" "Joe Patrik"
    Set-TextValue $sheet "B$row" "1996-02-13"
    Set-TextValue $sheet "C$row" "6935478966"
    Set-TextValue $sheet "D$row" "joepatrik@gmail.com"
    Set-TextValue $sheet "E$row" "Rajasthan"
    Set-TextValue $sheet "F$row" "Jaipur"
    Set-TextValue $sheet "G$row" "370081"
    Set-TextValue $sheet "H$row" "Pratik Colony, Near Yash Complex _x000d_`nSG Highway, Jaipur, Rajasthan"
    Set-TextValue $sheet "I$row" "1"
    Set-TextValue $sheet "J$row" "123456"
    Set-TextValue $sheet "K$row" "active"
}

function Set-RituAgarwalRow($sheet, $row) {
    Set-TextValue $sheet "A$row" "Ritu Agarwal"
    Set-TextValue $sheet "B$row" "2008-06-10"
    Set-TextValue $sheet "C$row" "54879672315"
    Set-TextValue $sheet "D$row" "rituagarwal@gmail.com"
    Set-TextValue $sheet "E$row" "Gujarat"
    Set-TextValue $sheet "F$row" "Gandhinagar"
    Set-TextValue $sheet "G$row" "345003"
    Set-TextValue $sheet "H$row" "Pratik Colony, Near Yash Complex _x000d_`nSG Highway, Jaipur, Rajasthan"
    Set-TextValue $sheet "I$row" "2"
    Set-TextValue $sheet "J$row" "123456"
    Set-TextValue $sheet "K$row" "pending"
}

Set-JoePatrikRow $memberSheet 3
Set-JoePatrikRow $memberSheet 5
Set-RituAgarwalRow $memberSheet 6
Set-JoePatrikRow $memberSheet 7
Set-RituAgarwalRow $memberSheet 8

# ---------------------------------------------------------------------------
# 4. New sheet booking_contact_details_tbl (appended at the end)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bookingSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$bookingSheet.Name = "booking_contact_details_tbl"

Set-TextValue $bookingSheet "A1" "name"
Set-TextValue $bookingSheet "B1" "email"
Set-TextValue $bookingSheet "C1" "phone"
Set-TextValue $bookingSheet "D1" "total_days"
Set-TextValue $bookingSheet "E1" "total_member"
Set-TextValue $bookingSheet "F1" "special_request"

Set-TextValue $bookingSheet "A3" "hetvee_sakaria"
Set-TextValue $bookingSheet "B3" "Hetu@2001"
Set-TextValue $bookingSheet "C3" "Hetvee Sakaria"

Set-TextValue $bookingSheet "A4" "hetvee_sakaria"
Set-TextValue $bookingSheet "B4" "Hetu@2001"
Set-TextValue $bookingSheet "C4" "Hetvee Sakaria"

# ---------------------------------------------------------------------------
# 5. Defined names -- update existing ranges & add the new one
# ---------------------------------------------------------------------------
$wb.Names.Item("admin_login_tbl").RefersTo = "=admin_login_tbl!`$A`$1:`$C`$5"
$wb.Names.Item("category_management_tbl").RefersTo = "=category_management_tbl!`$A`$1:`$B`$12"
$wb.Names.Item("member_signup_tbl").RefersTo = "=member_signup_tbl!`$A`$1:`$K`$8"
$wb.Names.Add("booking_contact_details_tbl", "=booking_contact_details_tbl!`$A`$1:`$F`$2")
